# Generate Report for Handoff
#
# The localization-status report is refreshed: the "3fd41955-..." row
# (row 4 on the Overview / zh-cn / de-de sheets) receives new
# "handoff generated" timestamps reflecting a fresh report generation.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-06 10:02:33"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-06 10:02:20"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-06 10:02:33"
